$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$pic = $s.Shapes.Item(4)
$pic.ActionSettings.Item(1).Hyperlink.Address = "https://example.com/"
